# Update erdri CDS definition sheet:
#  - D5 (3.1. Patient's status -> data_types) gets the new allowed-value list
#  - D8 (5.1. Age at onset -> data_types) text is replaced
#  - D9 (5.2. Age at diagnosis -> data_types) gets the same allowed-value list as D8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D8").Value = "antenatal, at birth, dd/mm/yyyy, undetermined"
$ws.Range("D9").Value = "antenatal, at birth, dd/mm/yyyy, undetermined"
$ws.Range("D5").Value = "alive, dead, lost in follow-up, opted-out"

$ws.Range("D12").Select()
